$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1024496205, "Eliana", "Lozano Romero", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001371417"),
    @(1032363637, "Camilo Antonio", "Monroy Peña", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001240633"),
    @(1032389108, "Luis Miguel", "Serrano Bermudez", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001353347"),
    @(1032419449, "Gerson Andres", "Rodriguez Rodriguez", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001421474"),
    @(1057581307, "Edwin Alexander", "Ramirez Perez", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001479088"),
    @(1082104140, "Oscar Libardo", "Lombana Charfuelan", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001413927"),
    @(1121859744, "Alfonso", "Martinez Andrade", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000013192"),
    @(1136881212, "German Andres", "Cabrera Rojas", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001468331"),
    @(14894713, "Jaime Leon", "Aguilar Arias", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000266299"),
    @(15813899, "Luis Carlos", "Belalcazar Ceron", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000249360"),
    @(15920734, "Jose Herney", "Ramirez Franco", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000773093"),
    @(16737949, "Rodrigo", "Jimenez Pizarro", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000177725"),
    @(19136869, "Marcelo Enrique", "Riveros Rojas", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001447216"),
    @(19187789, "Armando", "Espinosa Hernandez", "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000105244")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
